$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '29.852.75'
$ws.Range('E2').Value = '  -0.15%  '
$ws.Range('D3').Value = '1.887.29'
$ws.Range('E3').Value = '  -0.47%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.001'
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  +0.09%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '0.7464'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -4.98%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '242.39'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -0.60%  '
$ws.Range('E7').Value = '  +0.04%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.3114'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -0.91%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '25.29'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -2.07%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.07119'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -2.20%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.08491'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +4.74%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.7599'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -2.34%  '
$ws.Range('D13').Value = '1.895.80'
$ws.Range('E13').Value = '  -0.65%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '5.352'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -2.55%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '93.33'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -1.05%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '6.144'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -1.25%  '
$ws.Range('D17').Value = '29.919.74'
$ws.Range('E17').Value = '  +0.12%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '13.69'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  -2.06%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '243.26'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -1.55%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '0.000007793'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -0.49%  '
$ws.Range('D21').Value = '2.159.83'
$ws.Range('E21').Value = '  +0.66%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '1.001'
$ws.Range('D22').Style = 'Normal'
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '7.990'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -1.88%  '
$ws.Range('E24').Value = '  +0.08%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '0.1595'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +0.00%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '9.364'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -1.10%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '162.61'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -0.92%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '18.74'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -0.26%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '2.026'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -0.05%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '1.500'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +4.08%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '1.532'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -0.81%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '4.475'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -0.15%  '
$ws.Range('E33').Value = '  +0.50%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.05387'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -3.34%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.237'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -0.70%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.7443'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -1.54%  '
$ws.Range('E37').Value = '  -0.15%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '2.709'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +1.16%  '
$ws.Range('E39').Value = '  -0.22%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '2.767'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -0.87%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.4451'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -0.41%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '6.054'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +1.46%  '
$ws.Range('B43').Value = 'Maker'
$ws.Range('C43').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D43').Value = '1.091.23'
$ws.Range('E43').Value = '  -4.44%  '
$ws.Range('B44').Value = 'Aave'
$ws.Range('C44').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '72.50'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -2.27%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.8559'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +0.20%  '
$ws.Range('E46').Value = '  +0.02%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '102.26'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -0.02%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '7.675'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +1.57%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '1.859'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -2.15%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '3.065'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -2.45%  '
$ws.Range('D51').Value = '2.055.90'
$ws.Range('E51').Value = '  -0.06%  '
